$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the value for cell P2 to "Yoan" (reuses the existing shared string "Yoan")
$ws.Range("P2").Value = "Yoan"

# Move the active selection to P3, matching the post-edit cursor position
$ws.Range("P3").Select()
